$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Style = $ws.Range("H1").Style

$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 6
